$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 15.89417457580566
$ws.Range("D2").Value = 175

$ws.Range("C3").Value = 15.33699035644531
$ws.Range("D3").Value = 176

$ws.Range("C4").Value = 15.64908027648926
$ws.Range("D4").Value = 174

$ws.Range("C5").Value = 15.87915420532227
$ws.Range("D5").Value = 123

$ws.Range("C6").Value = 15.65909385681152
$ws.Range("D6").Value = 123
